$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update date strings in column A (rows 3-21): replace "/" with "-"
# NumberFormat is temporarily forced to Text ("@") so Excel does not
# auto-convert the dash-separated date text into a date serial number,
# then the cell's original style is restored so no visible style change
# is left behind.
for ($r = 3; $r -le 21; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $old = $cell.Value()
    $new = $old -replace '/', '-'
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $new
    $cell.Style = $origStyle
}

# Update attendance numbers that changed in the diff
$ws.Cells.Item(3, 4).Value = 1   # D3
$ws.Cells.Item(3, 7).Value = 1   # G3

$ws.Cells.Item(4, 4).Value = 1   # D4
$ws.Cells.Item(4, 5).Value = 1   # E4
$ws.Cells.Item(4, 8).Value = 0   # H4

$ws.Cells.Item(5, 4).Value = 1   # D5
$ws.Cells.Item(5, 5).Value = 1   # E5
$ws.Cells.Item(5, 8).Value = 0   # H5

$ws.Cells.Item(10, 4).Value = 1  # D10
$ws.Cells.Item(10, 5).Value = 1  # E10
$ws.Cells.Item(10, 8).Value = 0  # H10
